$d = $word.ActiveDocument

# --- Paragraph 1 (Heading1 title) ---
$d.Paragraphs(1).Range.Text = '"From Yellow Delight to Heavenly Bite: The Epic Journey of Banana Bread in the World of Food"'

# --- Paragraph 2 (intro paragraph) ---
$d.Paragraphs(2).Range.Text = "Banana bread is a beloved classic in the realm of comfort food. Its warm, moist texture and rich banana flavor make it a go-to treat for any occasion. This delightful baked good has stood the test of time and is a staple in many households. Whether enjoyed as a breakfast option or a sweet indulgence, banana bread has the power to evoke feelings of nostalgia and satisfaction. In this paper, we will delve into the fascinating world of food, using banana bread as a lens to explore the broader topic. By examining the origins, ingredients, and cultural significance of this delectable creation, we will uncover the intricate relationship between food and human experience. Join us on this tantalizing journey as we peel back the layers of banana bread's story, exploring the culinary delights and cultural implications that it brings to the table."

# --- Paragraph 3 (body paragraphs with manual line breaks) ---
$seg1 = "Banana bread is a classic and beloved treat that has been enjoyed by many for generations. Its sweet and moist texture, combined with the rich flavors of ripe bananas and warm spices, make it a popular choice for breakfast, snack, or even dessert. This delectable bread not only satisfies our cravings but also provides several health benefits."
$seg2 = "Firstly, banana bread is a great source of nutrients. Bananas themselves are loaded with essential vitamins and minerals, such as potassium, vitamin C, and vitamin B6. These nutrients help support our overall health and well-being. When incorporated into a bread recipe, bananas retain their nutritional value, making banana bread a delicious way to incorporate these vital nutrients into our diet."
$seg3 = "In addition to being nutritious, banana bread is also a great option for those with dietary restrictions. It can easily be made vegan or gluten-free by substituting ingredients such as eggs, butter, and flour with plant-based alternatives. This versatility allows individuals with specific dietary needs to still enjoy the comforting flavors of banana bread without compromising their health or lifestyle choices."
$seg4 = "Furthermore, banana bread is a great way to reduce food waste. Often, ripe bananas are overlooked and discarded when they become overripe and brown. However, these bananas are perfect for making banana bread. By transforming these seemingly unwanted bananas into a sweet and delicious treat, we can significantly reduce food waste and contribute to a more sustainable food system."
$seg5 = "Beyond the nutritional and sustainable aspects, banana bread also holds a special place in many people's hearts as a nostalgic comfort food. The warm and inviting aroma that fills the kitchen as the bread bakes brings back memories of home and family. Whether it's a treasured family recipe passed down through generations or a new discovery made in a cozy coffee shop, banana bread has the power to evoke feelings of comfort and happiness."
$seg6 = "In conclusion, banana bread is more than just a tasty treat; it is a versatile and nutritious food option that appeals to a wide range of individuals. Its ability to incorporate essential nutrients, accommodate dietary restrictions, reduce food waste, and evoke feelings of comfort and nostalgia make it a beloved favorite. So, the next time you have a few overripe bananas on your kitchen counter, consider whipping up a batch of banana bread and indulge in its delectable flavors while enjoying its many benefits."

$para3Text = $seg1 + "`v`v" + $seg2 + "`v`v" + $seg3 + "`v`v" + $seg4 + "`v`v" + $seg5 + "`v`v`v" + $seg6
$d.Paragraphs(3).Range.Text = $para3Text

# --- Paragraph 4 (conclusion paragraph) ---
$d.Paragraphs(4).Range.Text = "In conclusion, banana bread is a delightful and versatile treat that has gained popularity around the world. Its origins can be traced back to the Great Depression, where clever homemakers found a way to make use of overripe bananas. Today, banana bread stands as a testament to the ingenuity and resourcefulness of home cooks. With its moist, tender crumb and irresistible aroma, it has become a beloved comfort food for many. Whether enjoyed as a breakfast option, afternoon snack, or dessert, banana bread never fails to bring warmth and satisfaction. From classic recipes to innovative variations, there is a banana bread for every palate. So, next time you have a bunch of overripe bananas sitting on your counter, don't let them go to waste – whip up a batch of banana bread and indulge in its timeless deliciousness."

Write-Output "Done"
